# Splits the single run in a paragraph's text into one run per "word"
# (a maximal run of non-space characters) and one run per in-between
# run of whitespace -- matching how the target document represents
# e.g. "New guide and factsheet!" as 7 separate <w:r> elements:
# "New", " ", "guide", " ", "and", " ", "factsheet!".
#
# The Word object-model surface here has no direct "insert a run
# boundary" call, so we force one by toggling a character-formatting
# property (Bold) off/on/back-to-original on each segment after the
# first -- touching a Range's formatting makes the engine materialise
# it as its own run, and restoring the original value leaves the
# visible formatting unchanged.
function Split-ParagraphIntoWordRuns {
    param(
        $doc,
        $paraRange
    )

    $text = $paraRange.Text
    # A paragraph's Range.Text includes the trailing paragraph-mark
    # character(s); strip them before computing word/space offsets.
    $text = $text.TrimEnd([char]13, [char]10)
    if ($text.Length -eq 0) {
        return
    }

    $start = $paraRange.Start
    $matches = [System.Text.RegularExpressions.Regex]::Matches($text, '\S+|\s+')

    $segments = @()
    foreach ($m in $matches) {
        $segments += ,@($start + $m.Index, $start + $m.Index + $m.Length)
    }

    # Segment 0 needs no treatment: it is already everything left in
    # the original run once the later segments are carved off of it.
    # Walk the rest in order so each carve-off only ever touches the
    # still-unsplit tail of the run.
    for ($i = 1; $i -lt $segments.Count; $i++) {
        $segStart = $segments[$i][0]
        $segEnd = $segments[$i][1]
        $r = $doc.Range($segStart, $segEnd)
        $originalBold = $r.Bold
        $r.Bold = 1
        $r.Bold = $originalBold
    }
}

function Split-ParagraphByStyle {
    param(
        $doc,
        [string]$styleName
    )

    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Style.NameLocal -eq $styleName) {
            Split-ParagraphIntoWordRuns $doc $para.Range
            return
        }
    }
}

$d = $word.ActiveDocument

# "New guide and factsheet!" (Title paragraph)
Split-ParagraphByStyle $d "Title"

# "Tom Coleman" (Author paragraph)
Split-ParagraphByStyle $d "Author"

# "Guide on solving simultaneous equations available now!" (Abstract paragraph)
Split-ParagraphByStyle $d "Abstract"
